$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised values for existing rows 30-32 ---
$ws.Range("B30").Value = 69855
$ws.Range("J30").Value = 10510

$ws.Range("B31").Value = 71488
$ws.Range("J31").Value = 10705

$ws.Range("B32").Value = 74530
$ws.Range("I32").Value = 15833
$ws.Range("J32").Value = 11287

# --- Append new row 33 (Serie 01-08-2021) ---
# Format as text first so the date-like label is stored as a string
# (matching the other "Serie" labels in column A) instead of being
# auto-converted into a date serial number, then restore default style.
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "01-08-2021"
$ws.Range("A33").Style = "Normal"

$ws.Range("B33").Value = 80536
$ws.Range("C33").Value = 14475
$ws.Range("D33").Value = 9112
$ws.Range("E33").Value = 5447
$ws.Range("F33").Value = 5527
$ws.Range("G33").Value = 7250
$ws.Range("H33").Value = 11681
$ws.Range("I33").Value = 15490
$ws.Range("J33").Value = 11554
